$d = $word.ActiveDocument

# Locate the paragraph that currently holds the concatenated highlight
# color codes and split it into four separate paragraphs, each with its
# occurrence count appended.
$found = $d.Content.Find.Execute("#7cc867#fb5b89#f9cd59#c885da", $true, $false, $false, $false, $false,
                         $true, 1, $false, "#7cc867: 42^p#fb5b89: 18^p#f9cd59: 8^p#c885da: 8", 2)
